# إضافة حدث جديد في Card13 by HOSSAM at 2025-12-08 11:41:43
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# Fill in the previously-empty "nan" placeholder cells on the existing
# last row (row 20) for columns B..K and N.
$ws.Range("B20:K20").Value = "nan"
$ws.Range("N20").Value = "nan"

# Append the new service-log entry as row 21.
# Copy the card number from the row above so it keeps the same
# text (non-numeric) cell type as the rest of column A.
$ws.Range("A20").Copy($ws.Range("A21"))
# Columns B..K carry no reading for this entry - touch them (without
# writing any text) so the row keeps a full set of cells, matching the
# sheet's usual "blank placeholder" layout for unused columns.
$ws.Range("B21:K21").Borders.LineStyle = 0
$ws.Range("L21").Value = "24/3/2025"
$ws.Range("M21").Value = "تم تغير سير 1270"
$ws.Range("N21").Value = "قطع سير كويلر مسنن 1270"
$ws.Range("O21").Value = "فني"

$wb.Save()
